$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.226.11"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "1.645.11"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.74"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.93"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.872.79"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.29"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "1.681.07"
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "26.222.73"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.36"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.07"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.31"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -4.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.19"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.64"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0502"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("E34").Value = "  +1.67%  "
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").Value = "1.137.33"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.45"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D45").Value = "1.781.47"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.70"
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("E47").Value = "  +4.21%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.66"
$ws.Range("E50").Value = "  +2.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0969"
$ws.Range("E51").Value = "  +1.26%  "
